$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet right before the "总计" sheet, using the
#    existing "2021-Q4" sheet as a formatting template (same header/column
#    styling as every other quarter sheet).
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$ws.Name = "2022-Q1"

# NB: the sheet reference obtained *before* Worksheets.Add() tracks the
# workbook's sheet collection by position, so after the insert it now
# resolves to the freshly-added sheet instead of "总计" - re-resolve it by
# name afterwards to get a handle on the real summary sheet.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header-row formatting (B1:H1) and the column-A/index formatting
# (A2:H10) from the template sheet so the new sheet's style indices match
# the rest of the workbook.
$template.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2:H10").Copy()
$ws.Range("A2:H10").PasteSpecial(-4122)

# Columns B, D, E, F, G hold codes/percentages that must stay text (so
# leading zeros / trailing zeros are preserved) - force text format before
# writing the values.
$ws.Range("B2:B10").NumberFormat = "@"
$ws.Range("D2:G10").NumberFormat = "@"

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$data = @(
  @(0, "310358", "申万菱信新经济混合",                    "41.92", "77.32", "5.17", "2.1673", 3),
  @(1, "011488", "申万菱信乐享混合",                      "12.18", "77.88", "6.71", "0.8173", 2),
  @(2, "012210", "申万菱信智能汽车股票型证券投资基金A",       "4.76", "82.52", "6.80", "0.3237", 2),
  @(3, "012051", "申万菱信乐道三年持有期混合型证券投资基金",    "3.38", "81.64", "7.12", "0.2407", 2),
  @(4, "013634", "申万菱信双利混合A",                      "7.83", "22.26", "1.73", "0.1355", 5),
  @(5, "012211", "申万菱信智能汽车股票型证券投资基金C",       "1.40", "82.52", "6.80", "0.0952", 2),
  @(6, "001201", "申万菱信安鑫回报灵活配置混合A",            "4.72", "21.35", "0.79", "0.0373", 6),
  @(7, "001727", "申万菱信安鑫回报灵活配置混合C",            "3.96", "21.35", "0.79", "0.0313", 6),
  @(8, "013635", "申万菱信双利混合C",                      "0.75", "22.26", "1.73", "0.0130", 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $i + 2
  $d = $data[$i]
  $ws.Cells.Item($row, 1).Value = $d[0]
  $ws.Cells.Item($row, 2).Value = $d[1]
  $ws.Cells.Item($row, 3).Value = $d[2]
  $ws.Cells.Item($row, 4).Value = $d[3]
  $ws.Cells.Item($row, 5).Value = $d[4]
  $ws.Cells.Item($row, 6).Value = $d[5]
  $ws.Cells.Item($row, 7).Value = $d[6]
  $ws.Cells.Item($row, 8).Value = $d[7]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row right under the
#    header for the "2022-Q1" totals, pushing the older quarters down.
# ---------------------------------------------------------------------------
$totalSheet.Rows(2).Insert()

# The freshly inserted row picks up some stray formatting on B2:D2 - clear
# it, then restore the A-column index styling from the row below so it
# matches the rest of the table.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 9
$totalSheet.Range("D2").Value = 3.86

# Renumber the index column for the rows that got pushed down (they keep
# their original 0/1/2/3 values after the insert - bump them to 1/2/3/4).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
